$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Objetivos" long PT paragraph replaced by a professor name ---
$ws.Range("B10").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("C10").Value = '471420 - Carlos Antonio Reis Pereira Baptista'

# --- Fix up the style of brand-new B/C cells before filling them in: the sheet's
# <cols> defines overlapping ranges for column B (style 1 then style 2), so newly
# created B/C cells must have their number format/font/wrap copied explicitly from
# an existing, correctly-styled cell in the same column (B21/C21) rather than relying
# on the column default. ---
$ws.Range("B21").Copy() | Out-Null
$ws.Range("B19:B20").PasteSpecial(-4122)
$ws.Range("C21").Copy() | Out-Null
$ws.Range("C19:C20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 13-21: rewrite labels/values for the new compacted layout ---
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'
$ws.Range("C13").Value = '471420 - Carlos Antonio Reis Pereira Baptista'

$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '3480026 - João Paulo Pascon'
$ws.Range("C15").Value = '3480026 - João Paulo Pascon'

$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '7797767 - Viktor Pastoukhov'
$ws.Range("C18").Value = '7797767 - Viktor Pastoukhov'

$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'A avaliação será composta por duas provas (P1 e P2).'
$ws.Range("C19").Value = 'A avaliação será composta por duas provas (P1 e P2).'

$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total.'
$ws.Range("C20").Value = 'NS = NP1+NP2; NP1: questões da P1 valendo até 4p. no total; NP2: questões da P2 valendo até 6 p. no total.'

$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2.'
$ws.Range("C21").Value = 'A recuperação consistirá de uma prova de Recuperação (R), que irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2.'

# --- Row 22: Requisitos label only ---
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

# --- Row 23: requirement text (moved up from old row 26); trailing newline preserved ---
$ws.Range("A23").Clear()
$req23 = 'LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)' + "`n"
$ws.Range("B23").Value = $req23
$ws.Range("C23").Value = $req23

# --- Remove now-empty trailing rows 24-26 (old Bibliografia block + old Requisitos row) ---
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()

# --- Row heights to match new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
